$d = $word.ActiveDocument

function Get-ParaText($i) {
    return $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
}

function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ((Get-ParaText $i) -eq $text) {
            return $i
        }
    }
    return -1
}

function Rename-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Fix the "multiple polynomial equation system" (classification) list:
#    the databases generated from Eq. (1) (the "randP.../" prefixed items).
# ---------------------------------------------------------------------------

# Item 4: 100systems_100samplesPerSys.csv -> 100systems_10samplesPerSys.csv
$i = Find-ParaIndex("randPolynomialEquationSystem/100systems_100samplesPerSys.csv")
Rename-InParagraph $i "100systems_100samplesPerSys.csv" "100systems_10samplesPerSys.csv"

# Item 5: 100systems_1000samplesPerSys.csv -> 100systems_100samplesPerSys.csv
$i = Find-ParaIndex("randPolynomialEquationSystem/100systems_1000samplesPerSys.csv")
Rename-InParagraph $i "100systems_1000samplesPerSys.csv" "100systems_100samplesPerSys.csv"

# ---------------------------------------------------------------------------
# 2) Move the "And for the ones made from the Eq. (2)..." paragraph so that
#    it appears right after the item above (previously it appeared after the
#    6th item of this same list).
# ---------------------------------------------------------------------------

$afterIdx = Find-ParaIndex("randPolynomialEquationSystem/100systems_100samplesPerSys.csv")
$d.Paragraphs($afterIdx).Range.InsertParagraphAfter()
$newParaIdx = $afterIdx + 1
$newPara = $d.Paragraphs($newParaIdx)
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = $d.Styles("Normal")

$movedText = "And for the ones made from the Eq. (2), which were created with the same strategy (see databases in the directory databases/regression/multiplePolynomialEquationSystem), the following .csv (comma delimited) files were generated:"
$newPara.Range.Text = $movedText
$newPara.Range.Font.Name = "Times New Roman"
$newPara.Range.Font.Size = 12

# ---------------------------------------------------------------------------
# 3) Rename the 6th item of the first list (randP.../1000systems_1000samplesPerSys.csv)
#    to become p.../1systems_10samplesPerSys.csv
# ---------------------------------------------------------------------------

$i = Find-ParaIndex("randPolynomialEquationSystem/1000systems_1000samplesPerSys.csv")
Rename-InParagraph $i "randP" "p"
$i = Find-ParaIndex("polynomialEquationSystem/1000systems_1000samplesPerSys.csv")
Rename-InParagraph $i "1000systems_1000samplesPerSys.csv" "1systems_10samplesPerSys.csv"

# ---------------------------------------------------------------------------
# 4) Remove the original (now duplicate) "And for the ones..." paragraph.
# ---------------------------------------------------------------------------

$i = Find-ParaIndex($movedText)
# there will be two matches now (the moved one and the original) - remove the
# second occurrence (the original, later one).
$firstFound = $i
$secondFound = -1
for ($k = $firstFound + 1; $k -le $d.Paragraphs.Count; $k++) {
    if ((Get-ParaText $k) -eq $movedText) {
        $secondFound = $k
        break
    }
}
if ($secondFound -ne -1) {
    $d.Paragraphs($secondFound).Range.Delete()
}

# ---------------------------------------------------------------------------
# 5) Rename the 2nd, 3rd and 4th items of the second list (p.../ prefix).
# ---------------------------------------------------------------------------

$i = Find-ParaIndex("polynomialEquationSystem/1systems_10samplesPerSys.csv")
Rename-InParagraph $i "1systems_10samplesPerSys.csv" "10systems_10samplesPerSys.csv"

$i = Find-ParaIndex("polynomialEquationSystem/10systems_10samplesPerSys.csv")
Rename-InParagraph $i "10systems_10samplesPerSys.csv" "10systems_100samplesPerSys.csv"

$i = Find-ParaIndex("polynomialEquationSystem/10systems_100samplesPerSys.csv")
Rename-InParagraph $i "10systems_100samplesPerSys.csv" "100systems_10samplesPerSys.csv"

# ---------------------------------------------------------------------------
# 6) Delete the last two items of the second list entirely.
# ---------------------------------------------------------------------------

$i = Find-ParaIndex("polynomialEquationSystem/100systems_1000samplesPerSys.csv")
$d.Paragraphs($i).Range.Delete()

$i = Find-ParaIndex("polynomialEquationSystem/1000systems_1000samplesPerSys.csv")
$d.Paragraphs($i).Range.Delete()

# ---------------------------------------------------------------------------
# 7) Insert two new blank paragraphs right before the "Created in:" paragraph.
# ---------------------------------------------------------------------------

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = Get-ParaText $i
    if ($t.StartsWith("Created in:")) {
        break
    }
}
$createdInIdx = $i
$prev = $d.Paragraphs($createdInIdx - 1)
$prev.Range.InsertParagraphAfter()
$prev.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 8) Update the "Last update in:" date from November 26, 2021 to June 24, 2022.
# ---------------------------------------------------------------------------

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = Get-ParaText $i
    if ($t.StartsWith("Last update in:")) {
        break
    }
}
$lastUpdateIdx = $i
$p = $d.Paragraphs($lastUpdateIdx)
$rng = $p.Range
$rng.Find.Execute("November 26, 2021", $true, $false, $false, $false, $false, $true, 0, $false, "June 24, 2022", 2) | Out-Null

Write-Output "done"
